$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared-string text fixes ("fix typos and other gremlins") ---

# G3: "11.0; / 2.3, / 52" -> "... 52.4"
$ws.Range("G3").Value = "11.0; `n2.3, `n52.4"

# New cell I8 (previously empty) gets a new OR; 95% CI value
$ws.Range("I8").Value = "0.6;`n0.2,`n1.9"

# I12: "0.7; 0.3, 1.3" -> "1; 0.4, 2.3"
$ws.Range("I12").Value = "1;`n0.4,`n2.3"

# I13: "1.8; 1, 3.1" -> "3.3; 1.2, 10.1"
$ws.Range("I13").Value = "3.3;`n1.2,`n10.1"

# --- Numeric corrections ---
$ws.Range("O3").Value = 5
$ws.Range("K5").Value = 20
$ws.Range("O6").Value = 2
$ws.Range("N8").Value = 0
$ws.Range("O8").Value = 5
$ws.Range("L9").Value = 0
$ws.Range("O9").Value = 3
$ws.Range("K11").Value = 50
$ws.Range("M12").Value = 9
$ws.Range("N12").Value = 0
$ws.Range("O12").Value = 1
$ws.Range("M13").Value = 6
$ws.Range("N13").Value = 14
$ws.Range("O14").Value = 11
$ws.Range("M15").Value = 0
$ws.Range("O16").Value = 5

# --- Row 13 height normalized back to the sheet default ---
$ws.Rows(13).RowHeight = 35

# --- Column layout: narrow the raw-data columns away, widen the remaining ones ---
$ws.Columns("A").ColumnWidth = 20.493197278911566
$ws.Columns("B").ColumnWidth = 49.508503401360564
$ws.Range("C1:E1").EntireColumn.Hidden = $true
$ws.Columns("F").ColumnWidth = 11.794217687074866
$ws.Columns("G").ColumnWidth = 10.615646258503366
$ws.Columns("H").ColumnWidth = 15.799319727891165
$ws.Columns("I").ColumnWidth = 10.431972789115665
$ws.Columns("J").ColumnWidth = 10.064625850340166

# --- Selection moved to the last-edited cell ---
$ws.Range("J16").Select()
